$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Productdata": update SetupCosts (col E) and some StartingInventories
# (col C) values to reflect the new lot sizing rules.
# ---------------------------------------------------------------------------
$wsProd = $wb.Worksheets.Item("Productdata")

$wsProd.Range("E2").Value = 11.68992
$wsProd.Range("E3").Value = 2.171028
$wsProd.Range("E4").Value = 0.92261
$wsProd.Range("E5").Value = 1.43298
$wsProd.Range("E6").Value = 0.952776
$wsProd.Range("E7").Value = 0.2902439999999999
$wsProd.Range("E8").Value = 0.101376

$wsProd.Range("C9").Value = 928
$wsProd.Range("E9").Value = 0.8635200000000001

$wsProd.Range("C10").Value = 658
$wsProd.Range("E10").Value = 0.5153760000000001

$wsProd.Range("C11").Value = 1993
$wsProd.Range("E11").Value = 0.75747

$wsProd.Range("C12").Value = 457
$wsProd.Range("E12").Value = 1.235212

$wsProd.Range("C13").Value = 3736
$wsProd.Range("E13").Value = 12.57119999999999

$wsProd.Range("C14").Value = 1374
$wsProd.Range("E14").Value = 4.844166000000001

$wsProd.Range("C15").Value = 349
$wsProd.Range("E15").Value = 0.900614

$wsProd.Range("C16").Value = 488
$wsProd.Range("E16").Value = 0.881576

$wsProd.Range("C17").Value = 682
$wsProd.Range("E17").Value = 1.38996

$wsProd.Range("C18").Value = 206
$wsProd.Range("E18").Value = 0.45486

$wsProd.Range("C19").Value = 72
$wsProd.Range("E19").Value = 0.140288

$wsProd.Range("E20").Value = 63.10422000000001
$wsProd.Range("E21").Value = 66.991248
$wsProd.Range("E22").Value = 82.923264
$wsProd.Range("E23").Value = 255.152196

# ---------------------------------------------------------------------------
# Sheet "ForecastedAverageDemand": refresh forecast numbers per bucket.
# ---------------------------------------------------------------------------
$wsFAD = $wb.Worksheets.Item("ForecastedAverageDemand")

$wsFAD.Range("C2").Value = 600
$wsFAD.Range("D2").Value = 159
$wsFAD.Range("F2").Value = 298
$wsFAD.Range("G2").Value = 93
$wsFAD.Range("H2").Value = 37
$wsFAD.Range("I2").Value = 421
$wsFAD.Range("J2").Value = 297
$wsFAD.Range("K2").Value = 895
$wsFAD.Range("L2").Value = 207

$wsFAD.Range("C3").Value = 599
$wsFAD.Range("D3").Value = 149
$wsFAD.Range("F3").Value = 298
$wsFAD.Range("G3").Value = 89
$wsFAD.Range("H3").Value = 28
$wsFAD.Range("I3").Value = 419
$wsFAD.Range("J3").Value = 299
$wsFAD.Range("K3").Value = 908
$wsFAD.Range("L3").Value = 207

$wsFAD.Range("C4").Value = 599
$wsFAD.Range("D4").Value = 153
$wsFAD.Range("F4").Value = 297
$wsFAD.Range("G4").Value = 88
$wsFAD.Range("H4").Value = 29
$wsFAD.Range("I4").Value = 421
$wsFAD.Range("J4").Value = 302
$wsFAD.Range("K4").Value = 906
$wsFAD.Range("L4").Value = 212

$wsFAD.Range("C5").Value = 593
$wsFAD.Range("D5").Value = 150
$wsFAD.Range("F5").Value = 295
$wsFAD.Range("G5").Value = 91
$wsFAD.Range("H5").Value = 34
$wsFAD.Range("I5").Value = 419
$wsFAD.Range("J5").Value = 295
$wsFAD.Range("K5").Value = 898
$wsFAD.Range("L5").Value = 212

# ---------------------------------------------------------------------------
# Sheet "ForcastedStandardDeviation": refresh std-dev forecast numbers.
# ---------------------------------------------------------------------------
$wsFSD = $wb.Worksheets.Item("ForcastedStandardDeviation")

$wsFSD.Range("C2").Value = 75
$wsFSD.Range("D2").Value = 19.875
$wsFSD.Range("F2").Value = 37.25
$wsFSD.Range("G2").Value = 11.625
$wsFSD.Range("H2").Value = 4.625
$wsFSD.Range("I2").Value = 52.625
$wsFSD.Range("J2").Value = 37.125
$wsFSD.Range("K2").Value = 111.875
$wsFSD.Range("L2").Value = 25.875

$wsFSD.Range("C3").Value = 112.3125
$wsFSD.Range("D3").Value = 27.9375
$wsFSD.Range("F3").Value = 55.875
$wsFSD.Range("G3").Value = 16.6875
$wsFSD.Range("H3").Value = 5.25
$wsFSD.Range("I3").Value = 78.5625
$wsFSD.Range("J3").Value = 56.0625
$wsFSD.Range("K3").Value = 170.25
$wsFSD.Range("L3").Value = 38.8125

$wsFSD.Range("C4").Value = 131.03125
$wsFSD.Range("D4").Value = 33.46875
$wsFSD.Range("F4").Value = 64.96875
$wsFSD.Range("G4").Value = 19.25
$wsFSD.Range("H4").Value = 6.34375
$wsFSD.Range("I4").Value = 92.09375
$wsFSD.Range("J4").Value = 66.0625
$wsFSD.Range("K4").Value = 198.1875
$wsFSD.Range("L4").Value = 46.375

$wsFSD.Range("C5").Value = 138.984375
$wsFSD.Range("D5").Value = 35.15625
$wsFSD.Range("F5").Value = 69.140625
$wsFSD.Range("G5").Value = 21.328125
$wsFSD.Range("H5").Value = 7.96875
$wsFSD.Range("I5").Value = 98.203125
$wsFSD.Range("J5").Value = 69.140625
$wsFSD.Range("K5").Value = 210.46875
$wsFSD.Range("L5").Value = 49.6875

# ---------------------------------------------------------------------------
# Sheet "Capacity": update total capacity per bucket.
# ---------------------------------------------------------------------------
$wsCap = $wb.Worksheets.Item("Capacity")

$wsCap.Range("B2").Value = 50277.5
$wsCap.Range("B3").Value = 524868.75
$wsCap.Range("B4").Value = 524868.75
$wsCap.Range("B5").Value = 4948762.5
